$wb = $excel.ActiveWorkbook

# --- ALC row 15 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 833.3
$ws.Range("I15").Value = 833.3
$ws.Range("K15").Value = 2499.9
$ws.Range("M15").Value = -2330.9

# --- ALC row 96 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 125729.25
$ws.Range("I96").Value = 833.4286
$ws.Range("J96").Value = 1000000
$ws.Range("K96").Value = 2500.2858
$ws.Range("L96").Value = 3000000
$ws.Range("M96").Value = -1127.2858
$ws.Range("N96").Value = -3002746

# --- ALC row 98 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2684.2
$ws.Range("I98").Value = 2363
$ws.Range("K98").Value = 2363
$ws.Range("M98").Value = -865

# --- ALC row 122 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2684.2
$ws.Range("I122").Value = 2363
$ws.Range("K122").Value = 7089
$ws.Range("M122").Value = -4639

# --- ALC row 137 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3045.257
$ws.Range("I137").Value = 2768.3447
$ws.Range("J137").Value = 4383.6665
$ws.Range("K137").Value = 8305.0341
$ws.Range("L137").Value = 13150.9995
$ws.Range("M137").Value = -5755.034100000001
$ws.Range("N137").Value = -18250.9995

# --- ALC row 138 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3543.73
$ws.Range("J138").Value = 3889.872
$ws.Range("L138").Value = 11669.616
$ws.Range("N138").Value = -21949.616

# --- ARM row 32 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1544.76
$ws.Range("I32").Value = 1544.76
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1544.76
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1257.76
$ws.Range("N32").ClearContents()

# --- ARM row 45 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2802.1052
$ws.Range("J45").Value = 2861.625
$ws.Range("L45").Value = 2861.625
$ws.Range("N45").Value = -3615.625

# --- ARM row 80 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 69959.164
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 69959.164
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 69959.164
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -71955.164

# --- ARM row 83 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 69959.164
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 69959.164
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 209877.492
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -219861.492

# --- ARM row 102 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 6672653
$ws.Range("I102").Value = 3061.2307
$ws.Range("K102").Value = 3061.2307
$ws.Range("M102").Value = -1439.2307

# --- BSM row 99 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3589.3333
$ws.Range("I99").Value = 2070.25
$ws.Range("K99").Value = 2070.25
$ws.Range("M99").Value = -572.25

# --- BSM row 107 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2360
$ws.Range("I107").Value = 1624
$ws.Range("J107").Value = 4200
$ws.Range("K107").Value = 1624
$ws.Range("L107").Value = 4200
$ws.Range("M107").Value = 296
$ws.Range("N107").Value = -8040

# --- CRP row 58 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1491.1904
$ws.Range("I58").Value = 1258.4667
$ws.Range("J58").Value = 2073
$ws.Range("K58").Value = 1258.4667
$ws.Range("L58").Value = 2073
$ws.Range("M58").Value = -1055.4667
$ws.Range("N58").Value = -2479

# --- CRP row 99 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5815
$ws.Range("J99").Value = 6281.25
$ws.Range("L99").Value = 6281.25
$ws.Range("N99").Value = -9277.25

# --- CRP row 105 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 973.4
$ws.Range("I105").Value = 1023.7778
$ws.Range("K105").Value = 1023.7778
$ws.Range("M105").Value = 723.2222

# --- CRP row 126 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 5815
$ws.Range("J126").Value = 6281.25
$ws.Range("L126").Value = 18843.75
$ws.Range("N126").Value = -23783.75

# --- CRP row 132 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2501884.8
$ws.Range("I132").Value = 2858939.8
$ws.Range("K132").Value = 8576819.399999999
$ws.Range("M132").Value = -8574289.399999999

# --- CRP row 134 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2590.4
$ws.Range("I134").Value = 920.2727
$ws.Range("K134").Value = 2760.8181
$ws.Range("M134").Value = -225.8181

# --- CRP row 136 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1491.1904
$ws.Range("I136").Value = 1258.4667
$ws.Range("J136").Value = 2073
$ws.Range("K136").Value = 3775.4001
$ws.Range("L136").Value = 6219
$ws.Range("M136").Value = -1225.4001
$ws.Range("N136").Value = -11319

# --- CUL row 4 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4840051
$ws.Range("I4").Value = 2953875.2
$ws.Range("K4").Value = 8861625.600000001
$ws.Range("M4").Value = -8861513.600000001

# --- GSM row 52 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 43339.332
$ws.Range("J52").Value = 43339.332
$ws.Range("L52").Value = 43339.332
$ws.Range("N52").Value = -43857.332

# --- GSM row 102 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 18478.484
$ws.Range("I102").Value = 2213.6155
$ws.Range("J102").Value = 103055.8
$ws.Range("K102").Value = 2213.6155
$ws.Range("L102").Value = 103055.8
$ws.Range("M102").Value = -591.6154999999999
$ws.Range("N102").Value = -106299.8

# --- GSM row 126 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4312.5557
$ws.Range("I126").Value = 3785.5715
$ws.Range("K126").Value = 11356.7145
$ws.Range("M126").Value = -8886.7145

# --- LTW row 7 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7726.5557
$ws.Range("I7").Value = 6929.4614
$ws.Range("J7").Value = 9799
$ws.Range("K7").Value = 6929.4614
$ws.Range("L7").Value = 9799
$ws.Range("M7").Value = -6817.4614
$ws.Range("N7").Value = -10023

# --- LTW row 40 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7415.1904
$ws.Range("I40").Value = 6662.143
$ws.Range("K40").Value = 6662.143
$ws.Range("M40").Value = -6526.143

# --- LTW row 74 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 24230.691
$ws.Range("I74").Value = 14999
$ws.Range("K74").Value = 14999
$ws.Range("M74").Value = -14001

# --- LTW row 77 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H77").Value = 24230.691
$ws.Range("I77").Value = 14999
$ws.Range("K77").Value = 44997
$ws.Range("M77").Value = -40005

# --- LTW row 100 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 5159.8696
$ws.Range("I100").Value = 5613.357
$ws.Range("K100").Value = 5613.357
$ws.Range("M100").Value = -5072.357

# --- LTW row 126 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 7726.5557
$ws.Range("I126").Value = 6929.4614
$ws.Range("J126").Value = 9799
$ws.Range("K126").Value = 20788.3842
$ws.Range("L126").Value = 29397
$ws.Range("M126").Value = -18318.3842
$ws.Range("N126").Value = -34337

# --- LTW row 132 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3615.795
$ws.Range("I132").Value = 3518.4285
$ws.Range("J132").Value = 3729.389
$ws.Range("K132").Value = 10555.2855
$ws.Range("L132").Value = 11188.167
$ws.Range("M132").Value = -8025.2855
$ws.Range("N132").Value = -16248.167

# --- WVR row 17 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 4000
$ws.Range("I17").Value = 4000
$ws.Range("K17").Value = 4000
$ws.Range("M17").Value = -3828

# --- WVR row 45 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 21812.334
$ws.Range("J45").Value = 21812.334
$ws.Range("L45").Value = 21812.334
$ws.Range("N45").Value = -22794.334

# --- WVR row 126 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3203.6843
$ws.Range("I126").Value = 3483.1667
$ws.Range("K126").Value = 10449.5001
$ws.Range("M126").Value = -7979.500100000001

# --- WVR row 131 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H131").Value = 95000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 95000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 95000
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -105080

# --- WVR row 132 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1525.8276
$ws.Range("I132").Value = 1384.5769
$ws.Range("K132").Value = 4153.7307
$ws.Range("M132").Value = -1623.7307
